$d = $word.ActiveDocument

# These custom character styles each have only a Bold run property set
# (besides Color); the wml.xsd schema requires <w:b/> before <w:color/>
# inside <w:rPr>, but they were authored with <w:color/> first. Re-assign
# Font.Bold to itself so the style's rPr gets re-serialized in the
# schema-correct order (b, color).
foreach ($name in @("KeywordTok","ImportTok","ControlFlowTok","AlertTok","ErrorTok")) {
    $s = $d.Styles($name)
    $s.Font.Bold = $s.Font.Bold
}

# Same issue, but for styles with only Italic set (besides Color):
# re-order to (i, color).
foreach ($name in @("CommentTok","DocumentationTok")) {
    $s = $d.Styles($name)
    $s.Font.Italic = $s.Font.Italic
}

# Styles with both Bold and Italic set (besides Color): re-order to
# (b, i, color).
foreach ($name in @("AnnotationTok","CommentVarTok","InformationTok","WarningTok")) {
    $s = $d.Styles($name)
    $s.Font.Bold = $s.Font.Bold
    $s.Font.Italic = $s.Font.Italic
}
